$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A243").Value = 0
$ws.Range("B243").Value = 5.205180650487272
$ws.Range("A244").Value = 1
$ws.Range("B244").Value = 3.956576556616381
$ws.Range("A245").Value = 2
$ws.Range("B245").Value = 3.473479429979991
$ws.Range("A246").Value = 3
$ws.Range("B246").Value = 3.437121623022951
$ws.Range("A247").Value = 4
$ws.Range("B247").Value = 2.643228509312538
$ws.Range("A248").Value = 5
$ws.Range("B248").Value = 1.760924730453898
$ws.Range("A249").Value = 6
$ws.Range("B249").Value = 1.523583777472727
$ws.Range("A250").Value = 7
$ws.Range("B250").Value = 1.357597039083963
$ws.Range("A251").Value = 8
$ws.Range("B251").Value = 1.306134976117688
$ws.Range("A252").Value = 9
$ws.Range("B252").Value = 1.266161055231813
$ws.Range("A253").Value = 10
$ws.Range("B253").Value = 1.201616142328877
$ws.Range("A254").Value = 11
$ws.Range("B254").Value = 1.195585881234915
$ws.Range("A255").Value = 12
$ws.Range("B255").Value = 1.226689229916379
$ws.Range("A256").Value = 13
$ws.Range("B256").Value = 1.180434916728189
$ws.Range("A257").Value = 14
$ws.Range("B257").Value = 1.115839579263878
$ws.Range("A258").Value = 15
$ws.Range("B258").Value = 1.052252726258405
$ws.Range("A259").Value = 16
$ws.Range("B259").Value = 1.087150177696135
$ws.Range("A260").Value = 17
$ws.Range("B260").Value = 1.178073317154063
$ws.Range("A261").Value = 18
$ws.Range("B261").Value = 1.099288381657827
$ws.Range("A262").Value = 19
$ws.Range("B262").Value = 1.058448205402254
$ws.Range("A263").Value = 20
$ws.Range("B263").Value = 1.080840272589291
$ws.Range("A264").Value = 21
$ws.Range("B264").Value = 1.080944215594006
$ws.Range("A265").Value = 22
$ws.Range("B265").Value = 1.096330180071861
$ws.Range("A266").Value = 23
$ws.Range("B266").Value = 1.087868660421187
$ws.Range("A267").Value = 24
$ws.Range("B267").Value = 1.009688846451093
$ws.Range("A268").Value = 25
$ws.Range("B268").Value = 1.072081390473391
$ws.Range("A269").Value = 26
$ws.Range("B269").Value = 1.079625517010108
$ws.Range("A270").Value = 27
$ws.Range("B270").Value = 1.094251279083718
$ws.Range("A271").Value = 28
$ws.Range("B271").Value = 1.138035285145209
$ws.Range("A272").Value = 29
$ws.Range("B272").Value = 1.076502317283534
$ws.Range("A273").Value = 30
$ws.Range("B273").Value = 1.089079628275054
$ws.Range("A274").Value = 31
$ws.Range("B274").Value = 1.097446472860466
$ws.Range("A275").Value = 32
$ws.Range("B275").Value = 1.088708718929113
$ws.Range("A276").Value = 33
$ws.Range("B276").Value = 1.091054144781229
$ws.Range("A277").Value = 34
$ws.Range("B277").Value = 1.135269469533876
$ws.Range("A278").Value = 35
$ws.Range("B278").Value = 1.054888496761779
$ws.Range("A279").Value = 36
$ws.Range("B279").Value = 0.9804809349498624
$ws.Range("A280").Value = 37
$ws.Range("B280").Value = 0.9913985729168245
$ws.Range("A281").Value = 38
$ws.Range("B281").Value = 1.021971301165635
$ws.Range("A282").Value = 39
$ws.Range("B282").Value = 1.023071655321763
$ws.Range("A283").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B283").Value = 1.023071655321763
$ws.Range("A284").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B284").Value = 0.9613599579775841
$ws.Range("A285").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B285").Value = 0.9892343708303535
$ws.Range("A286").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B286").Value = 0.9549625858457806
$ws.Range("A287").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B287").Value = 0.9627892913692859
$ws.Range("A288").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B288").Value = 0.9401874004788022
$ws.Range("A289").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B289").Value = 0.938167371880488
$ws.Range("A290").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B290").Value = 1.05464406564183
$ws.Range("A291").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B291").Value = 0.9791052752838884
$ws.Range("A292").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B292").Value = 0.9294519271927117
$ws.Range("A293").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B293").Value = 0.9558654946082695
$ws.Range("A294").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B294").Value = 1.00506258823306
$ws.Range("A295").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B295").Value = 0.9948401599271234
$ws.Range("A296").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B296").Value = 0.9343099875625508
$ws.Range("A297").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B297").Value = 1.001491826711858
$ws.Range("A298").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B298").Value = 0.9559397363245918
$ws.Range("A299").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B299").Value = 1.024585566394081
$ws.Range("A300").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B300").Value = 1.037604621111003
$ws.Range("A301").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B301").Value = 0.9585248935763495
$ws.Range("A302").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B302").Value = 1.072542206164056
$ws.Range("A303").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B303").Value = 1.034233152250161
$ws.Range("A304").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B304").Value = 1.010083856371503
$ws.Range("A305").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B305").Value = 0.9784211237516599
$ws.Range("A306").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B306").Value = 0.9485652003365996
$ws.Range("A307").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B307").Value = 0.9956109935629377
$ws.Range("A308").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B308").Value = 0.9532199319171614
$ws.Range("A309").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B309").Value = 0.940111663437014
$ws.Range("A310").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B310").Value = 0.9906518691492431
$ws.Range("A311").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B311").Value = 1.016246213042413
$ws.Range("A312").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B312").Value = 1.01033331917733
$ws.Range("A313").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B313").Value = 0.9564526746337625
$ws.Range("A314").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B314").Value = 0.9380522297971179
$ws.Range("A315").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B315").Value = 0.9644038827912014
$ws.Range("A316").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B316").Value = 1.005472059523607
$ws.Range("A317").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B317").Value = 0.9692668236645543
$ws.Range("A318").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B318").Value = 0.8782583283145834
$ws.Range("A319").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B319").Value = 0.9811203487043575
$ws.Range("A320").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B320").Value = 0.9398007544985979
$ws.Range("A321").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B321").Value = 1.007863459260986
$ws.Range("A322").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B322").Value = 0.9752874356258755
$ws.Range("A323").Value = '<__main__.DisplayOutputs object at 0x7f5c780c0940>'
$ws.Range("B323").Value = 0.9996396006527586
